$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename "price" -> "description" and "image" -> "image_url"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "image_url"

# Data row 2: replace scraped Amazon data with scraped Jumia data
$ws.Range("A2").Value = "https://www.jumia.com.eg/ar/catalog/?q=blouse"

# B2 becomes an empty cell (kept present, but with no content).
# Touching a formatting property after clearing the value keeps the
# (now-empty) cell node alive in the saved worksheet XML.
$ws.Range("B2").Value = ""
$ws.Range("B2").Font.Bold = $false

$ws.Range("C2").Value = "شميز حريمي أبيض ساده فورمال236.00 جنيه295.00 جنيه20%4.3 out of 5(167)"
$ws.Range("D2").Value = "data:image/gif;base64,R0lGODlhAQABAIAAAAAAAP///yH5BAEAAAAALAAAAAABAAEAAAIBRAA7"
